$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "49.669.92"
Set-TextValue "E2" "  +3.21%  "

Set-TextValue "D3" "2.614.62"
Set-TextValue "E3" "  +4.54%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.09%  "

Set-TextValue "B5" "BNB"
Set-TextValue "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D5" "323.74"
Set-TextValue "E5" "  +0.75%  "

Set-TextValue "B6" "Solana"
Set-TextValue "C6" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D6" "109.93"
Set-TextValue "E6" "  +1.68%  "

Set-TextValue "D7" "0.534"
Set-TextValue "E7" "  +1.56%  "

Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.07%  "

Set-TextValue "D9" "0.560"
Set-TextValue "E9" "  +3.73%  "

Set-TextValue "D10" "40.83"
Set-TextValue "E10" "  +2.47%  "

Set-TextValue "D11" "20.58"
Set-TextValue "E11" "  +1.99%  "

Set-TextValue "D12" "0.0823"
Set-TextValue "E12" "  +0.91%  "

Set-TextValue "E13" "  +0.65%  "

Set-TextValue "D14" "7.32"
Set-TextValue "E14" "  +2.15%  "

Set-TextValue "D15" "3.030.77"
Set-TextValue "E15" "  +4.87%  "

Set-TextValue "D16" "2.606.82"
Set-TextValue "E16" "  +4.11%  "

Set-TextValue "D17" "0.871"
Set-TextValue "E17" "  +3.23%  "

Set-TextValue "D18" "49.699.70"
Set-TextValue "E18" "  +3.62%  "

Set-TextValue "D19" "3.10"
Set-TextValue "E19" "  +12.11%  "

Set-TextValue "D20" "13.34"
Set-TextValue "E20" "  +2.09%  "

Set-TextValue "D21" "6.77"
Set-TextValue "E21" "  +0.12%  "

Set-TextValue "D22" "0.0₃0953"
Set-TextValue "E22" "  +0.85%  "

Set-TextValue "D23" "281.07"
Set-TextValue "E23" "  +1.81%  "

Set-TextValue "D24" "72.78"
Set-TextValue "E24" "  +1.23%  "

Set-TextValue "D25" "2.57"
Set-TextValue "E25" "  +0.99%  "

Set-TextValue "D26" "26.61"
Set-TextValue "E26" "  +3.31%  "

Set-TextValue "E27" "  +0.01%  "

Set-TextValue "E28" "  -5.06%  "

Set-TextValue "D29" "9.96"
Set-TextValue "E29" "  +1.68%  "

Set-TextValue "D30" "0.144"
Set-TextValue "E30" "  +3.33%  "

Set-TextValue "D31" "36.27"
Set-TextValue "E31" "  +2.73%  "

Set-TextValue "D32" "49.60"
Set-TextValue "E32" "  +0.89%  "

Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "5.46"
Set-TextValue "E33" "  +2.45%  "

Set-TextValue "B34" "Celestia"
Set-TextValue "C34" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D34" "19.62"
Set-TextValue "E34" "  +0.42%  "

Set-TextValue "D35" "1.01"
Set-TextValue "E35" "  +0.17%  "

Set-TextValue "D36" "0.0795"
Set-TextValue "E36" "  +1.36%  "

Set-TextValue "D37" "2.05"
Set-TextValue "E37" "  +5.10%  "

Set-TextValue "D38" "4.78"

Set-TextValue "D39" "3.08"
Set-TextValue "E39" "  +5.71%  "

Set-TextValue "D40" "22.88"
Set-TextValue "E40" "  +6.95%  "

Set-TextValue "B41" "Stellar"
Set-TextValue "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.112"
Set-TextValue "E41" "  +0.85%  "

Set-TextValue "B42" "Monero"
Set-TextValue "C42" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D42" "123.30"
Set-TextValue "E42" "  +1.81%  "

Set-TextValue "E43" "  +0.34%  "

Set-TextValue "D44" "0.0315"
Set-TextValue "E44" "  +3.91%  "

Set-TextValue "E45" "  +6.36%  "

Set-TextValue "D46" "2.051.53"
Set-TextValue "E46" "  +2.55%  "

Set-TextValue "E47" "  +10.58%  "

Set-TextValue "E48" "  +9.62%  "

Set-TextValue "D49" "9.04"
Set-TextValue "E49" "  +1.10%  "

Set-TextValue "D50" "5.37"
Set-TextValue "E50" "  +4.00%  "

Set-TextValue "D51" "81.99"
Set-TextValue "E51" "  +2.71%  "

Write-Host "Cryptos list updated"